$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44195
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 13500
$ws.Range("N2").Value = "`$/caja 15 kilos"
$ws.Range("P2").Value = 900
$ws.Range("Q2").Value = 15

$ws.Range("D3").Value = 44349
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 31000
$ws.Range("N3").Value = "`$/malla 25 kilos"
$ws.Range("P3").Value = 1240
$ws.Range("Q3").Value = 25

$ws.Range("D4").Value = 44230
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 24500
$ws.Range("N4").Value = "`$/malla 25 kilos"
$ws.Range("P4").Value = 980
$ws.Range("Q4").Value = 25

$ws.Range("D5").Value = 44321
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 22000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 23000
$ws.Range("N5").Value = "`$/malla 25 kilos"
$ws.Range("P5").Value = 920
$ws.Range("Q5").Value = 25

$ws.Range("D6").Value = 44265
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 23000
$ws.Range("M6").Value = 21500
$ws.Range("N6").Value = "`$/malla 25 kilos"
$ws.Range("P6").Value = 860
$ws.Range("Q6").Value = 25

$ws.Range("D7").Value = 44258
$ws.Range("J7").Value = 1600
$ws.Range("K7").Value = 26000
$ws.Range("L7").Value = 28000
$ws.Range("M7").Value = 27000
$ws.Range("N7").Value = "`$/malla 25 kilos"
$ws.Range("P7").Value = 1080
$ws.Range("Q7").Value = 25

$ws.Range("D8").Value = 44203
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 19000
$ws.Range("M8").Value = 18500
$ws.Range("N8").Value = "`$/malla 25 kilos"
$ws.Range("P8").Value = 740
$ws.Range("Q8").Value = 25

$ws.Range("D9").Value = 44231
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 24500
$ws.Range("N9").Value = "`$/malla 25 kilos"
$ws.Range("P9").Value = 980
$ws.Range("Q9").Value = 25

$ws.Range("D10").Value = 44252
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 21000
$ws.Range("L10").Value = 23000
$ws.Range("M10").Value = 22000
$ws.Range("N10").Value = "`$/malla 25 kilos"
$ws.Range("P10").Value = 880
$ws.Range("Q10").Value = 25

$ws.Range("D11").Value = 44237
$ws.Range("J11").Value = 1100
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 21000
$ws.Range("M11").Value = 20500
$ws.Range("N11").Value = "`$/malla 25 kilos"
$ws.Range("P11").Value = 820
$ws.Range("Q11").Value = 25

$ws.Range("D12").Value = 44223
$ws.Range("J12").Value = 1700
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 22000
$ws.Range("M12").Value = 21000
$ws.Range("N12").Value = "`$/malla 25 kilos"
$ws.Range("P12").Value = 840
$ws.Range("Q12").Value = 25

$ws.Range("D13").Value = 44209
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14000
$ws.Range("N13").Value = "`$/caja 15 kilos"
$ws.Range("P13").Value = 933
$ws.Range("Q13").Value = 15

$ws.Range("D14").Value = 44322
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 23000
$ws.Range("N14").Value = "`$/malla 25 kilos"
$ws.Range("P14").Value = 920
$ws.Range("Q14").Value = 25

$ws.Range("D15").Value = 44300
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 23000
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = 24000
$ws.Range("N15").Value = "`$/malla 25 kilos"
$ws.Range("P15").Value = 960
$ws.Range("Q15").Value = 25

$ws.Range("D16").Value = 44224
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 20000
$ws.Range("L16").Value = 22000
$ws.Range("M16").Value = 21000
$ws.Range("N16").Value = "`$/malla 25 kilos"
$ws.Range("P16").Value = 840
$ws.Range("Q16").Value = 25

$ws.Range("D17").Value = 44259
$ws.Range("J17").Value = 600
$ws.Range("K17").Value = 25000
$ws.Range("L17").Value = 26000
$ws.Range("M17").Value = 25500
$ws.Range("N17").Value = "`$/malla 25 kilos"
$ws.Range("P17").Value = 1020
$ws.Range("Q17").Value = 25

$ws.Range("D18").Value = 44251
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 21000
$ws.Range("L18").Value = 23000
$ws.Range("M18").Value = 22000
$ws.Range("N18").Value = "`$/malla 25 kilos"
$ws.Range("P18").Value = 880
$ws.Range("Q18").Value = 25

$ws.Range("D19").Value = 44245
$ws.Range("J19").Value = 1100
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 23000
$ws.Range("M19").Value = 21500
$ws.Range("N19").Value = "`$/malla 25 kilos"
$ws.Range("P19").Value = 860
$ws.Range("Q19").Value = 25

$ws.Range("D20").Value = 44210
$ws.Range("J20").Value = 700
$ws.Range("K20").Value = 23000
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = 24000
$ws.Range("N20").Value = "`$/malla 25 kilos"
$ws.Range("P20").Value = 960
$ws.Range("Q20").Value = 25

$ws.Range("D21").Value = 44266
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = 21500
$ws.Range("N21").Value = "`$/malla 25 kilos"
$ws.Range("P21").Value = 860
$ws.Range("Q21").Value = 25

$ws.Range("D22").Value = 44238
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 21000
$ws.Range("L22").Value = 23000
$ws.Range("M22").Value = 22000
$ws.Range("N22").Value = "`$/malla 25 kilos"
$ws.Range("P22").Value = 880
$ws.Range("Q22").Value = 25

$ws.Range("D23").Value = 44202
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 19000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 19500
$ws.Range("N23").Value = "`$/malla 25 kilos"
$ws.Range("P23").Value = 780
$ws.Range("Q23").Value = 25

$ws.Range("D24").Value = 44188
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14500
$ws.Range("N24").Value = "`$/caja 15 kilos"
$ws.Range("P24").Value = 967
$ws.Range("Q24").Value = 15

$ws.Range("D25").Value = 44216
$ws.Range("J25").Value = 800
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 15500
$ws.Range("N25").Value = "`$/caja 15 kilos"
$ws.Range("P25").Value = 1033
$ws.Range("Q25").Value = 15

$ws.Range("D26").Value = 44189
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14500
$ws.Range("N26").Value = "`$/caja 15 kilos"
$ws.Range("P26").Value = 967
$ws.Range("Q26").Value = 15

$ws.Range("D27").Value = 44314
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 23000
$ws.Range("L27").Value = 25000
$ws.Range("M27").Value = 24000
$ws.Range("N27").Value = "`$/malla 25 kilos"
$ws.Range("P27").Value = 960
$ws.Range("Q27").Value = 25

$ws.Range("D28").Value = 44350
$ws.Range("J28").Value = 140
$ws.Range("K28").Value = 28000
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = 29000
$ws.Range("N28").Value = "`$/malla 25 kilos"
$ws.Range("P28").Value = 1160
$ws.Range("Q28").Value = 25

$ws.Range("D29").Value = 44272
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 21000
$ws.Range("L29").Value = 23000
$ws.Range("M29").Value = 22000
$ws.Range("N29").Value = "`$/malla 25 kilos"
$ws.Range("P29").Value = 880
$ws.Range("Q29").Value = 25

$ws.Range("D30").Value = 44294
$ws.Range("J30").Value = 400
$ws.Range("K30").Value = 20000
$ws.Range("L30").Value = 23000
$ws.Range("M30").Value = 21500
$ws.Range("N30").Value = "`$/malla 25 kilos"
$ws.Range("P30").Value = 860
$ws.Range("Q30").Value = 25

$ws.Range("D31").Value = 44308
$ws.Range("J31").Value = 360
$ws.Range("K31").Value = 23000
$ws.Range("L31").Value = 25000
$ws.Range("M31").Value = 24000
$ws.Range("N31").Value = "`$/malla 25 kilos"
$ws.Range("P31").Value = 960
$ws.Range("Q31").Value = 25

$ws.Range("D32").Value = 44315
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 23000
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = 24000
$ws.Range("N32").Value = "`$/malla 25 kilos"
$ws.Range("P32").Value = 960
$ws.Range("Q32").Value = 25

$ws.Range("D33").Value = 44329
$ws.Range("J33").Value = 600
$ws.Range("K33").Value = 25000
$ws.Range("L33").Value = 27000
$ws.Range("M33").Value = 26000
$ws.Range("N33").Value = "`$/malla 25 kilos"
$ws.Range("P33").Value = 1040
$ws.Range("Q33").Value = 25

$ws.Range("D34").Value = 44279
$ws.Range("J34").Value = 700
$ws.Range("K34").Value = 24000
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = 24500
$ws.Range("N34").Value = "`$/malla 25 kilos"
$ws.Range("P34").Value = 980
$ws.Range("Q34").Value = 25

$ws.Range("D35").Value = 44280
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 24000
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = 24500
$ws.Range("N35").Value = "`$/malla 25 kilos"
$ws.Range("P35").Value = 980
$ws.Range("Q35").Value = 25

$ws.Range("D36").Value = 44293
$ws.Range("J36").Value = 700
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 23000
$ws.Range("M36").Value = 21500
$ws.Range("N36").Value = "`$/malla 25 kilos"
$ws.Range("P36").Value = 860
$ws.Range("Q36").Value = 25

$ws.Range("D37").Value = 44343
$ws.Range("J37").Value = 700
$ws.Range("K37").Value = 25000
$ws.Range("L37").Value = 27000
$ws.Range("M37").Value = 26000
$ws.Range("N37").Value = "`$/malla 25 kilos"
$ws.Range("P37").Value = 1040
$ws.Range("Q37").Value = 25

$ws.Range("D38").Value = 44175
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 20000
$ws.Range("L38").Value = 21000
$ws.Range("M38").Value = 20500
$ws.Range("N38").Value = "`$/caja 15 kilos"
$ws.Range("P38").Value = 1367
$ws.Range("Q38").Value = 15

$ws.Range("D39").Value = 44286
$ws.Range("J39").Value = 1200
$ws.Range("K39").Value = 20000
$ws.Range("L39").Value = 23000
$ws.Range("M39").Value = 21500
$ws.Range("N39").Value = "`$/malla 25 kilos"
$ws.Range("P39").Value = 860
$ws.Range("Q39").Value = 25

$ws.Range("D40").Value = 44181
$ws.Range("J40").Value = 900
$ws.Range("K40").Value = 12000
$ws.Range("L40").Value = 13000
$ws.Range("M40").Value = 12500
$ws.Range("N40").Value = "`$/caja 15 kilos"
$ws.Range("P40").Value = 833
$ws.Range("Q40").Value = 15

$ws.Range("D41").Value = 44328
$ws.Range("J41").Value = 500
$ws.Range("K41").Value = 25000
$ws.Range("L41").Value = 27000
$ws.Range("M41").Value = 26000
$ws.Range("N41").Value = "`$/malla 25 kilos"
$ws.Range("P41").Value = 1040
$ws.Range("Q41").Value = 25

$ws.Range("D42").Value = 44301
$ws.Range("J42").Value = 1000
$ws.Range("K42").Value = 23000
$ws.Range("L42").Value = 25000
$ws.Range("M42").Value = 24000
$ws.Range("N42").Value = "`$/malla 25 kilos"
$ws.Range("P42").Value = 960
$ws.Range("Q42").Value = 25

$ws.Range("D43").Value = 44217
$ws.Range("J43").Value = 760
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 16000
$ws.Range("M43").Value = 15500
$ws.Range("N43").Value = "`$/caja 15 kilos"
$ws.Range("P43").Value = 1033
$ws.Range("Q43").Value = 15

$ws.Range("D44").Value = 44244
$ws.Range("J44").Value = 1700
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 23000
$ws.Range("M44").Value = 21500
$ws.Range("N44").Value = "`$/malla 25 kilos"
$ws.Range("P44").Value = 860
$ws.Range("Q44").Value = 25

$ws.Range("D45").Value = 44307
$ws.Range("J45").Value = 600
$ws.Range("K45").Value = 23000
$ws.Range("L45").Value = 25000
$ws.Range("M45").Value = 24000
$ws.Range("N45").Value = "`$/malla 25 kilos"
$ws.Range("P45").Value = 960
$ws.Range("Q45").Value = 25

$ws.Range("D46").Value = 44273
$ws.Range("J46").Value = 300
$ws.Range("K46").Value = 20000
$ws.Range("L46").Value = 22000
$ws.Range("M46").Value = 21000
$ws.Range("N46").Value = "`$/malla 25 kilos"
$ws.Range("P46").Value = 840
$ws.Range("Q46").Value = 25
